# Daily attendance processing - 2025-10-25 20:21:45
# Normalizes the "Recorded By" (column G) entries so that when the
# literal entry "System" is present in the comma-separated list of
# recorders, it is moved to the front of that list (the relative order
# of the remaining entries is preserved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }
    if ($raw.GetType().Name -ne "String") { continue }

    $parts = $raw.Split(",")
    if ($parts.Length -lt 2) { continue }

    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $lastIndex = $trimmed.Length - 1
    if ($trimmed[$lastIndex].Equals("System")) {
        $rest = $trimmed[0..($lastIndex - 1)]
        $newParts = @("System") + $rest
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value = $newValue
    }
}
